$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.407.12'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('E3').Value = '  -1.84%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'588.20"
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').Value = "'176.77"
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value = "'0.611"
$ws.Range('E7').Value = '  +1.64%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = "'1.00"
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '3.457.07'
$ws.Range('E9').Value = '  -1.80%  '
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('D11').Value = "'6.97"
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('D12').Value = "'0.418"
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('D13').Value = '4.055.59'
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').Value = "'29.58"
$ws.Range('E15').Value = '  -3.64%  '
$ws.Range('D16').Value = '66.282.25'
$ws.Range('E16').Value = '  -0.95%  '
$ws.Range('D17').Value = "'0.0000172"
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('D18').Value = '3.461.44'
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('E19').Value = '  -2.40%  '
$ws.Range('E20').Value = '  -2.02%  '
$ws.Range('D21').Value = "'374.26"
$ws.Range('E21').Value = '  -2.05%  '
$ws.Range('E22').Value = '  -3.32%  '
$ws.Range('D23').Value = "'73.38"
$ws.Range('E23').Value = '  +2.34%  '
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = "'0.537"
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').Value = "'0.0000126"
$ws.Range('E26').Value = '  +3.39%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '3.606.03'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = "'9.91"
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = "'0.178"
$ws.Range('E29').Value = '  +2.49%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = "'0.999"
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = "'5.86"
$ws.Range('E31').Value = '  -2.58%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = "'2.01"
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = "'23.69"
$ws.Range('E33').Value = '  -3.89%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Value = "'1.00"
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = "'1.28"
$ws.Range('E35').Value = '  -7.14%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = "'7.04"
$ws.Range('E36').Value = '  -3.15%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = "'1.55"
$ws.Range('E37').Value = '  -1.39%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = "'161.42"
$ws.Range('E38').Value = '  +1.62%  '
$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D39').Value = "'0.884"
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Value = "'28.52"
$ws.Range('E40').Value = '  -2.87%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = "'1.82"
$ws.Range('E41').Value = '  +0.69%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = "'2.64"
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.770.74'
$ws.Range('E44').Value = '  +1.67%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = "'6.45"
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = "'0.0693"
$ws.Range('E46').Value = '  -2.11%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = "'25.35"
$ws.Range('E47').Value = '  -1.49%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = "'337.95"
$ws.Range('E48').Value = '  +2.79%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = "'40.02"
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = "'0.0293"
$ws.Range('E50').Value = '  -2.23%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = "'0.103"
$ws.Range('E51').Value = '  +0.26%  '
